$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.732.60"
$ws.Range('E2').Value = '  +0.43%  '
$ws.Range('D3').Value = "'1.853.09"
$ws.Range('E3').Value = '  +0.55%  '
$ws.Range('D4').Value = "'1.003"
$ws.Range('E4').Value = '  +0.30%  '
$ws.Range('D5').Value = "'312.68"
$ws.Range('E5').Value = '  -0.62%  '
$ws.Range('D6').Value = "'1.003"
$ws.Range('E6').Value = '  +0.24%  '
$ws.Range('D7').Value = "'0.4281"
$ws.Range('E7').Value = '  +1.23%  '
$ws.Range('D8').Value = "'0.3590"
$ws.Range('E8').Value = '  -1.29%  '
$ws.Range('D9').Value = "'0.07305"
$ws.Range('E9').Value = '  +0.22%  '
$ws.Range('D10').Value = "'0.8751"
$ws.Range('E10').Value = '  -1.50%  '
$ws.Range('E11').Value = '  +0.41%  '
$ws.Range('D12').Value = "'1.812.41"
$ws.Range('E12').Value = '  -0.57%  '
$ws.Range('D13').Value = "'6.551"
$ws.Range('E13').Value = '  -0.20%  '
$ws.Range('D14').Value = "'5.334"
$ws.Range('E14').Value = '  +0.05%  '
$ws.Range('D15').Value = "'0.07007"
$ws.Range('E15').Value = '  +1.87%  '
$ws.Range('D16').Value = "'1.006"
$ws.Range('E16').Value = '  +0.45%  '
$ws.Range('D17').Value = "'79.64"
$ws.Range('E17').Value = '  +0.65%  '
$ws.Range('D18').Value = "'0.000008957"
$ws.Range('E18').Value = '  +1.03%  '
$ws.Range('E19').Value = '  +0.40%  '
$ws.Range('E20').Value = '  -0.96%  '
$ws.Range('D21').Value = "'27.796.80"
$ws.Range('E21').Value = '  +0.67%  '
$ws.Range('D22').Value = "'5.001"
$ws.Range('E22').Value = '  +0.29%  '
$ws.Range('E23').Value = '  -1.65%  '
$ws.Range('D24').Value = "'2.097.82"
$ws.Range('E24').Value = '  +1.39%  '
$ws.Range('D25').Value = "'1.990"
$ws.Range('E25').Value = '  +4.58%  '
$ws.Range('D26').Value = "'155.07"
$ws.Range('E26').Value = '  +0.94%  '
$ws.Range('D27').Value = "'18.54"
$ws.Range('E27').Value = '  -2.81%  '
$ws.Range('D28').Value = "'120.62"
$ws.Range('E28').Value = '  -1.99%  '
$ws.Range('D29').Value = "'5.265"
$ws.Range('E29').Value = '  -0.51%  '
$ws.Range('D30').Value = "'1.883"
$ws.Range('E30').Value = '  -0.31%  '
$ws.Range('D31').Value = "'0.08922"
$ws.Range('E31').Value = '  -0.18%  '
$ws.Range('D32').Value = "'0.7583"
$ws.Range('E32').Value = '  -1.40%  '
$ws.Range('D33').Value = "'2.968"
$ws.Range('E33').Value = '  +0.86%  '
$ws.Range('E34').Value = '  -1.12%  '
$ws.Range('D35').Value = "'1.125"
$ws.Range('E35').Value = '  +2.63%  '
$ws.Range('D36').Value = "'1.003"
$ws.Range('E36').Value = '  +0.25%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = "'1.109"
$ws.Range('E37').Value = '  +0.60%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = "'0.05439"
$ws.Range('E38').Value = '  +1.15%  '
$ws.Range('D39').Value = "'0.01933"
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('E40').Value = '  +0.57%  '
$ws.Range('D41').Value = "'0.1672"
$ws.Range('E41').Value = '  +1.10%  '
$ws.Range('D42').Value = "'0.5095"
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('D43').Value = "'6.631"
$ws.Range('E43').Value = '  -3.53%  '
$ws.Range('D44').Value = "'8.424"
$ws.Range('E44').Value = '  +1.82%  '
$ws.Range('D45').Value = "'106.23"
$ws.Range('E45').Value = '  +1.83%  '
$ws.Range('D46').Value = "'0.06528"
$ws.Range('E46').Value = '  -0.97%  '
$ws.Range('D47').Value = "'10.34"
$ws.Range('E47').Value = '  -0.47%  '
$ws.Range('D48').Value = "'0.4679"
$ws.Range('E48').Value = '  -1.10%  '
$ws.Range('E49').Value = '  +0.26%  '
$ws.Range('E50').Value = '  -0.47%  '
$ws.Range('D51').Value = "'1.792"
$ws.Range('E51').Value = '  +2.25%  '

# Reset style on text-forced price cells so no stray quote-prefix style/format lingers
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
